# Update automatico via Actualizar 02-21-2021 12-24-31
#
# The "Disponibilidad" tracker gets a fresh batch of rows appended by the
# scheduled updater each run; the timestamp written into column D shifts
# down a block (14 rows) at a time as older batches age out one slot.
# This pass just refreshes column D with the newest captured timestamps -
# the most recent batch (rows 2-15) gets the brand new capture time, and
# the two older batches (rows 16-29, 30-43) shift to the next newest
# timestamps, exactly mirroring what the previous run wrote one slot up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newestTimestamp = 44248.51686204165
$middleTimestamp = 44248.49556476852
$oldestTimestamp = 44248.47426738426

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 4).Value = $newestTimestamp
}

for ($row = 16; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Value = $middleTimestamp
}

for ($row = 30; $row -le 43; $row++) {
    $ws.Cells.Item($row, 4).Value = $oldestTimestamp
}
